$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold, centered style) from the existing header cell E1
# onto the two new header cells F1:G1, then set their text.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "conf.low"
$ws.Range("G1").Value = "conf.high"

# Fill in the conf.low / conf.high numeric values for rows 2-12
$data = @(
    @(0.2417851131769609, 0.4041498271749452),
    @(-0.07242325075447811, 0.01487718907357428),
    @(-0.0765219875772896, 0.01575835661650615),
    @(-0.08554231145938371, 0.005664659719890885),
    @(-0.0845742750861021, 0.006376522622121683),
    @(0.1536543828898994, 0.2499125356388241),
    @(-0.00120729519242606, 0.001604776582844399),
    @(-0.07638550470031472, -0.0005340969787968192),
    @(-0.04320404511485494, 0.04027778528032111),
    @(-0.03900474074773055, 0.0449678237729194),
    @(-0.04737522532439117, 0.04705772082858172)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 6).Value = $data[$i][0]
    $ws.Cells.Item($r, 7).Value = $data[$i][1]
}

# Update the used range dimension to reflect the new columns
$ws.Range("A1:G12").Select() | Out-Null
